$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 81
$ws1.Range("F3").Value = 11856
$ws1.Range("F6").Value = 353
$ws1.Range("F8").Value = 11772
$ws1.Range("F10").Value = 1173
$ws1.Range("F12").Value = 51
$ws1.Range("F13").Value = 1776
$ws1.Range("F14").Value = 5834
$ws1.Range("F16").Value = 3533
$ws1.Range("F17").Value = 186

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 81
$ws4.Range("F5").Value = 11856
$ws4.Range("F9").Value = 353
$ws4.Range("F11").Value = 11772
$ws4.Range("F13").Value = 1173
$ws4.Range("F15").Value = 51
$ws4.Range("F16").Value = 1776
$ws4.Range("F18").Value = 5834
$ws4.Range("F20").Value = 3533
$ws4.Range("F21").Value = 186
